$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the resume summary text in C1 (Friday resume update)
$ws.Range("C1").Value = "Software quality analyst with 3.2 years of experience in software quality processes, Involved in end-to-end features testing. Skills in Automation testing, Manual UI Testing, Database and API Testing."

# Size column C to fit the (now much longer) summary text, as Excel's
# "AutoFit selection" does whenever a cell's content changes significantly.
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 172

# Move/restore the active selection to C2
$ws.Range("C2").Select() | Out-Null
